$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 239-244 (dates 44313-44318, i.e. 2021-04-27 .. 2021-05-02)
$data = @(
    @(239, 44313, 1, 9,  59.38242280285036),
    @(240, 44314, 1, 9,  59.38242280285036),
    @(241, 44315, 6, 14, 92.37265769332278),
    @(242, 44316, 1, 13, 85.77461071522829),
    @(243, 44317, 1, 14, 92.37265769332278),
    @(244, 44318, 3, 14, 92.37265769332278)
)

foreach ($row in $data) {
    $r = $row[0]
    # Copy the date-column style (column A) from the row above so the new
    # date cell keeps the same number format / style index.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
